{"js": "// Replace the multiplication-expression text in the worksheet table.\n// Each entry maps the original \"AxB=\" cell text to its replacement.\nconst replacements = [\n  [\"71\u00d725=\", \"43\u00d741=\"],\n  [\"35\u00d723=\", \"40\u00d745=\"],\n  [\"72\u00d736=\", \"82\u00d791=\"],\n  [\"43\u00d722=\", \"74\u00d773=\"],\n  [\"20\u00d782=\", \"29\u00d759=\"],\n  [\"78\u00d714=\", \"14\u00d761=\"],\n  [\"54\u00d712=\", \"53\u00d798=\"],\n  [\"84\u00d797=\", \"31\u00d728=\"],\n  [\"46\u00d769=\", \"95\u00d798=\"],\n  [\"96\u00d731=\", \"72\u00d768=\"],\n  [\"24\u00d720=\", \"84\u00d794=\"],\n  [\"60\u00d722=\", \"97\u00d764=\"],\n  [\"22\u00d722=\", \"62\u00d771=\"],\n  [\"62\u00d729=\", \"30\u00d778=\"],\n  [\"83\u00d757=\", \"11\u00d760=\"],\n  [\"66\u00d787=\", \"71\u00d753=\"],\n  [\"43\u00d785=\", \"89\u00d737=\"],\n  [\"84\u00d724=\", \"50\u00d771=\"],\n  [\"57\u00d779=\", \"79\u00d773=\"],\n  [\"31\u00d751=\", \"85\u00d717=\"],\n  [\"18\u00d789=\", \"65\u00d788=\"],\n  [\"49\u00d742=\", \"12\u00d794=\"],\n  [\"74\u00d714=\", \"80\u00d723=\"],\n  [\"50\u00d770=\", \"68\u00d748=\"],\n  [\"30\u00d774=\", \"60\u00d791=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-expression text in the worksheet table.\n# Each pair is the original \"AxB=\" cell text and its replacement.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"71\u00d725=\", \"43\u00d741=\"),\n    @(\"35\u00d723=\", \"40\u00d745=\"),\n    @(\"72\u00d736=\", \"82\u00d791=\"),\n    @(\"43\u00d722=\", \"74\u00d773=\"),\n    @(\"20\u00d782=\", \"29\u00d759=\"),\n    @(\"78\u00d714=\", \"14\u00d761=\"),\n    @(\"54\u00d712=\", \"53\u00d798=\"),\n    @(\"84\u00d797=\", \"31\u00d728=\"),\n    @(\"46\u00d769=\", \"95\u00d798=\"),\n    @(\"96\u00d731=\", \"72\u00d768=\"),\n    @(\"24\u00d720=\", \"84\u00d794=\"),\n    @(\"60\u00d722=\", \"97\u00d764=\"),\n    @(\"22\u00d722=\", \"62\u00d771=\"),\n    @(\"62\u00d729=\", \"30\u00d778=\"),\n    @(\"83\u00d757=\", \"11\u00d760=\"),\n    @(\"66\u00d787=\", \"71\u00d753=\"),\n    @(\"43\u00d785=\", \"89\u00d737=\"),\n    @(\"84\u00d724=\", \"50\u00d771=\"),\n    @(\"57\u00d779=\", \"79\u00d773=\"),\n    @(\"31\u00d751=\", \"85\u00d717=\"),\n    @(\"18\u00d789=\", \"65\u00d788=\"),\n    @(\"49\u00d742=\", \"12\u00d794=\"),\n    @(\"74\u00d714=\", \"80\u00d723=\"),\n    @(\"50\u00d770=\", \"68\u00d748=\"),\n    @(\"30\u00d774=\", \"60\u00d791=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
